# Applies the Config.xlsx update described in the commit:
# 1. Repoints the "FirstSheetName" column (C) for the three
#    "...ActionRequired_New_..." rows on DataFetchXL from
#    "Transmittals_New" to the new "Transmittals_New_ActionRequired" value.
# 2. Adds a new "FLD_Transmittals_LeftNavigationBar" test-row to both the
#    DataFetchFlag and DataFetchXL sheets (including a hyperlink on the
#    DataFetchXL sheet pointing at the new test-data workbook).

$wb = $excel.ActiveWorkbook

$wsFlag = $wb.Worksheets.Item("DataFetchFlag")
$wsXL   = $wb.Worksheets.Item("DataFetchXL")

# ---------------------------------------------------------------------
# 1. DataFetchXL: repoint FirstSheetName for the "ActionRequired_New_*" rows
# ---------------------------------------------------------------------
$wsXL.Range("C12").Value = "Transmittals_New_ActionRequired"
$wsXL.Range("C13").Value = "Transmittals_New_ActionRequired"
$wsXL.Range("C14").Value = "Transmittals_New_ActionRequired"

# ---------------------------------------------------------------------
# 2. DataFetchFlag: append row 18
# ---------------------------------------------------------------------
$wsFlag.Range("A18").Value = "FLD_Transmittals_LeftNavigationBar"
$wsFlag.Range("B18").Value = "XL"
$wsFlag.Range("B18").Style = "Normal"

# Keep the list validation covering the newly added row.
$wsFlag.Range("B2:B18").Validation.Delete()
$wsFlag.Range("B2:B18").Validation.Add(3, 1, 1, "XL,DB")
$wsFlag.Range("B2:B18").Validation.IgnoreBlank = $true
$wsFlag.Range("B2:B18").Validation.InCellDropdown = $true
$wsFlag.Range("B2:B18").Validation.ShowInput = $true
$wsFlag.Range("B2:B18").Validation.ShowError = $true

# ---------------------------------------------------------------------
# 3. DataFetchXL: append row 18 (with hyperlink on column B)
# ---------------------------------------------------------------------
$wsXL.Range("A18").Value = "FLD_Transmittals_LeftNavigationBar"

$target = "file:///\\src\com\proj\suiteTRANSMITTALS\testdata\Transmittals-Fluid-Navigation.xlsx"
$wsXL.Hyperlinks.Add($wsXL.Range("B18"), $target, "", "", "\\src\\com\\proj\\suiteTRANSMITTALS\\testdata\\Transmittals-Fluid-Navigation.xlsx")
$wsXL.Range("B18").Style = "Hyperlink"

$wsXL.Range("C18").Value = "Fluid_Navigation"

$wsFlag.Select()
$wsFlag.Range("B18").Select()
